$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need to be forced to text so they
# remain stored as strings (matching the original inline-string cell type)
# instead of being auto-converted to numbers by Excel.
$textForceCells = @("D5", "D6", "D9", "D10", "D11", "D17", "D19", "D22", "D23", "D24", "D25", "D28", "D30", "D32", "D35", "D37", "D38", "D42", "D43", "D45", "D47", "D49", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.698.35"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "1.632.66"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "213.33"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "0.500"
$ws.Range("E6").Value = "  +3.15%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").Value = "0.0623"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").Value = "19.24"
$ws.Range("E10").Value = "  +2.49%  "
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("D12").Value = "1.858.25"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "1.635.28"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "26.664.24"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").Value = "63.55"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "218.40"
$ws.Range("E19").Value = "  +8.05%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "6.16"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "9.37"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  +4.79%  "
$ws.Range("D25").Value = "147.74"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D28").Value = "6.86"
$ws.Range("E28").Value = "  +3.97%  "
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("D30").Value = "0.0505"
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "3.31"
$ws.Range("E32").Value = "  +3.83%  "
$ws.Range("E33").Value = "  +2.00%  "
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").Value = "2.39"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "1.220.06"
$ws.Range("E36").Value = "  +5.02%  "
$ws.Range("D37").Value = "0.0172"
$ws.Range("E37").Value = "  +5.21%  "
$ws.Range("D38").Value = "0.806"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").Value = "0.795"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").Value = "5.33"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").Value = "1.766.02"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").Value = "92.61"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("E46").Value = "  +2.61%  "
$ws.Range("D47").Value = "55.30"
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").Value = "0.0512"
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("D50").Value = "7.62"
$ws.Range("E50").Value = "  +4.20%  "
$ws.Range("E51").Value = "  -0.29%  "

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
